$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 13; this shifts existing rows 13-20 down to 14-21
$ws.Rows.Item(13).Insert()

# Populate the new row 13 with data (copy of the "Espárragos" record pattern, new date/volume/prices)
$ws.Range("A13").Value = 5
$ws.Range("B13").Value = "Macroferia Regional de Talca"
$ws.Range("C13").Value = "Maule"
$ws.Range("D13").Value = 44468
$ws.Range("E13").Value = 7
$ws.Range("F13").Value = 300000000
$ws.Range("G13").Value = "Espárragos"
$ws.Range("H13").Value = "Verde"
$ws.Range("I13").Value = "Primera"
$ws.Range("J13").Value = 3000
$ws.Range("K13").Value = 1500
$ws.Range("L13").Value = 1500
$ws.Range("M13").Value = 1500
$ws.Range("N13").Value = "$/kilo"
$ws.Range("O13").Value = "Provincia de Linares"
$ws.Range("P13").Value = 1500
$ws.Range("Q13").Value = 1
$ws.Range("R13").Value = "Hortaliza"
